$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in new test case rows (3, 4, 5) column-by-column, reflecting the
# additional login test scenarios (Invalid Password, Invalid Username,
# Empty Credentials) that now produce an "Error" / 401 response.
$ws.Range("A3").Value = "TC-LOGIN-002"
$ws.Range("B3").Value = "Invalid Password"
$ws.Range("A4").Value = "TC-LOGIN-003"
$ws.Range("B4").Value = "Invalid Username"
$ws.Range("A5").Value = "TC-LOGIN-004"
$ws.Range("B5").Value = "Empty Credentials"

$ws.Range("D3").Value = "testpwd"
$ws.Range("E3").Value = "Error"
$ws.Range("C4").Value = "test"

$ws.Range("C3").Value = "admin"
$ws.Range("D4").Value = "password"
$ws.Range("E4").Value = "Error"
$ws.Range("E5").Value = "Error"

$ws.Range("F3").Value = 401
$ws.Range("F4").Value = 401
$ws.Range("F5").Value = 401

# Widen column B to fit the new scenario text
$ws.Columns.Item(2).ColumnWidth = 15.333333

# Update selection to the last edited cell
$ws.Range("F5").Select()
